# Updated cryptos list on Mon Nov 13 23:24:34 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "36.584.69";    E = "  -1.35%  " },
    @{ Row = 3;  E = "  +0.86%  " },
    @{ Row = 4;  E = "  -0.17%  " },
    @{ Row = 5;  D = "241.85";       E = "  -2.35%  " },
    @{ Row = 6;  D = "0.665";        E = "  +0.27%  " },
    @{ Row = 8;  D = "52.67";        E = "  -6.70%  " },
    @{ Row = 9;  D = "58.93";        E = "  -1.65%  " },
    @{ Row = 10; D = "0.359";        E = "  -6.48%  " },
    @{ Row = 11; E = "  -3.72%  " },
    @{ Row = 12; E = "  -1.05%  " },
    @{ Row = 13; D = "0.896";        E = "  -0.18%  " },
    @{ Row = 14; D = "14.50";        E = "  -9.33%  " },
    @{ Row = 15; D = "2.365.94";     E = "  +0.38%  " },
    @{ Row = 16; E = "  -5.32%  " },
    @{ Row = 17; D = "2.037.78";     E = "  -0.67%  " },
    @{ Row = 18; D = "36.525.61";    E = "  -1.66%  " },
    @{ Row = 19; D = "16.34";        E = "  -14.30%  " },
    @{ Row = 20; D = "71.61";        E = "  -4.30%  " },
    @{ Row = 21; D = "0.0₃0863"; E = "  -3.14%  " },
    @{ Row = 22; D = "5.27";         E = "  -2.54%  " },
    @{ Row = 23; D = "235.63";       E = "  -0.50%  " },
    @{ Row = 24; E = "  +0.21%  " },
    @{ Row = 25; E = "  -4.49%  " },
    @{ Row = 26; D = "9.28";         E = "  -2.97%  " },
    @{ Row = 27; D = "2.12";         E = "  -2.49%  " },
    @{ Row = 28; D = "163.24";       E = "  -4.49%  " },
    @{ Row = 29; D = "20.54";        E = "  +2.40%  " },
    @{ Row = 30; E = "  -1.48%  " },
    @{ Row = 31; D = "5.09";         E = "  +0.78%  " },
    @{ Row = 32; D = "1.15";         E = "  -3.20%  " },
    @{ Row = 33; D = "4.55";         E = "  -1.57%  " },
    @{ Row = 34; E = "  -3.90%  " },
    @{ Row = 36; D = "2.30";         E = "  +1.82%  " },
    @{ Row = 37; E = "  -0.19%  " },
    @{ Row = 38; D = "0.0820";       E = "  -6.62%  " },
    @{ Row = 39; E = "  -6.11%  " },
    @{ Row = 40; D = "2.93";         E = "  -5.05%  " },
    @{ Row = 41; D = "4.84";         E = "  -5.58%  " },
    @{ Row = 42; E = "  -3.37%  " },
    @{ Row = 43; D = "1.12";         E = "  -2.44%  " },
    @{ Row = 44; E = "  -5.76%  " },
    @{ Row = 45; D = "93.87";        E = "  -2.93%  " },
    @{ Row = 46; D = "1.401.72";     E = "  +9.46%  " },
    @{ Row = 47; D = "15.61";        E = "  -9.42%  " },
    @{ Row = 48; D = "7.36";         E = "  +8.30%  " },
    @{ Row = 49; D = "2.32";         E = "  -3.15%  " },
    @{ Row = 50; E = "  +0.08%  " },
    @{ Row = 51; D = "2.254.51";     E = "  +0.55%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Force the cell to stay text. Several of these price strings
        # (e.g. "241.85") parse as valid numbers and Excel would silently
        # convert/reformat them, losing the exact textual representation
        # used by the source feed. Applying a text format first, then
        # resetting the style back to Normal afterwards, keeps the cell's
        # type as text without leaving a custom number format behind.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
